$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes/rotates the full data of certain rows among themselves
# (cyclic re-ordering of records), while row numbers, formatting and the
# Y/AA "Startdatum"/"Slutdatum" columns (identical "2026-02-07" in every
# row of the affected set) stay put. The mapping below gives, for each
# destination row, the row whose current ("before edit") content should
# end up there.
$srcForDst = @{
    6  = 7
    7  = 8
    8  = 6
    9  = 10
    10 = 9
    15 = 18
    16 = 17
    17 = 15
    18 = 16
    21 = 22
    22 = 23
    23 = 21
    32 = 33
    33 = 32
}

# All rows that participate (both as a source and a destination).
$rows = @(6,7,8,9,10,15,16,17,18,21,22,23,32,33)

# Snapshot each row's values BEFORE any writes happen, because several
# rows both feed into and receive data from other rows in the same batch
# (this is a set of rotation cycles, not independent pairwise swaps).
# Two column blocks are used (A:X and AB:AY) so column Y/Z/AA (Startdatum/
# Starttid/Slutdatum) is left completely untouched -- writing through
# Value/Value2 would otherwise let Excel "helpfully" reinterpret the
# "2026-02-07" text as a real date serial, which is not what the source
# file has (plain text cells) and would add spurious number formatting.
$snapLeft = @{}
$snapRight = @{}
foreach ($r in $rows) {
    $snapLeft[$r]  = $ws.Range("A${r}:X${r}").Value2()
    $snapRight[$r] = $ws.Range("AB${r}:AY${r}").Value2()
}

foreach ($r in $rows) {
    $s = $srcForDst[$r]
    $ws.Range("A${r}:X${r}").Value2   = $snapLeft[$s]
    $ws.Range("AB${r}:AY${r}").Value2 = $snapRight[$s]
}
